$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps its original plain-text representation
# (values like "71.394.97" or "0.190" must not be re-interpreted as numbers).
$ws.Range("D2:D51").NumberFormat = "@"

# Apply the updated cell values from the source diff
$ws.Range('D2').Value = '71.394.97'
$ws.Range('E2').Value = '  +0.72%  '
$ws.Range('D3').Value = '3.647.03'
$ws.Range('E3').Value = '  +0.40%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '582.47'
$ws.Range('E5').Value = '  -1.44%  '
$ws.Range('D6').Value = '189.54'
$ws.Range('E6').Value = '  -2.40%  '
$ws.Range('D7').Value = '3.642.48'
$ws.Range('E7').Value = '  +0.50%  '
$ws.Range('D8').Value = '0.629'
$ws.Range('E8').Value = '  -2.04%  '
$ws.Range('E9').Value = '  +0.06%  '
$ws.Range('D10').Value = '0.190'
$ws.Range('E10').Value = '  +3.56%  '
$ws.Range('D11').Value = '0.664'
$ws.Range('E11').Value = '  -1.29%  '
$ws.Range('D12').Value = '55.18'
$ws.Range('E12').Value = '  -4.87%  '
$ws.Range('D13').Value = '0.0000311'
$ws.Range('E13').Value = '  -0.12%  '
$ws.Range('D14').Value = '9.77'
$ws.Range('E14').Value = '  -1.32%  '
$ws.Range('D15').Value = '4.224.83'
$ws.Range('E15').Value = '  +0.16%  '
$ws.Range('D16').Value = '19.88'
$ws.Range('E16').Value = '  -2.14%  '
$ws.Range('D17').Value = '3.635.73'
$ws.Range('E17').Value = '  -0.05%  '
$ws.Range('D18').Value = '71.193.15'
$ws.Range('E18').Value = '  +0.47%  '
$ws.Range('D19').Value = '12.69'
$ws.Range('E19').Value = '  -0.45%  '
$ws.Range('E20').Value = '  -0.81%  '
$ws.Range('E21').Value = '  +1.24%  '
$ws.Range('D22').Value = '506.12'
$ws.Range('E22').Value = '  +3.35%  '
$ws.Range('E23').Value = '  +3.56%  '
$ws.Range('D24').Value = '5.01'
$ws.Range('E24').Value = '  -2.32%  '
$ws.Range('E25').Value = '  -1.10%  '
$ws.Range('D26').Value = '96.74'
$ws.Range('E26').Value = '  +6.56%  '
$ws.Range('D27').Value = '11.65'
$ws.Range('E27').Value = '  +2.11%  '
$ws.Range('D28').Value = '3.04'
$ws.Range('E28').Value = '  -3.59%  '
$ws.Range('D29').Value = '9.53'
$ws.Range('E29').Value = '  -0.01%  '
$ws.Range('D30').Value = '7.83'
$ws.Range('E30').Value = '  -2.39%  '
$ws.Range('E31').Value = '  -1.33%  '
$ws.Range('E32').Value = '  +3.81%  '
$ws.Range('D33').Value = '66.60'
$ws.Range('E33').Value = '  -0.53%  '
$ws.Range('D35').Value = '580.34'
$ws.Range('E35').Value = '  -5.03%  '
$ws.Range('D36').Value = '3.25'
$ws.Range('E36').Value = '  +8.99%  '
$ws.Range('D37').Value = '39.47'
$ws.Range('E37').Value = '  -2.65%  '
$ws.Range('D38').Value = '0.416'
$ws.Range('E38').Value = '  +1.29%  '
$ws.Range('D39').Value = '0.0₃0812'
$ws.Range('E39').Value = '  -4.03%  '
$ws.Range('D40').Value = '0.998'
$ws.Range('E40').Value = '  -0.07%  '
$ws.Range('D41').Value = '3.35'
$ws.Range('E41').Value = '  +5.10%  '
$ws.Range('E42').Value = '  -2.32%  '
$ws.Range('D43').Value = '0.139'
$ws.Range('E44').Value = '  -1.52%  '
$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D45').Value = '0.0459'
$ws.Range('E45').Value = '  +0.02%  '
$ws.Range('B46').Value = 'ApeXProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D46').Value = '3.56'
$ws.Range('E46').Value = '  +5.25%  '
$ws.Range('D47').Value = '3.253.19'
$ws.Range('E47').Value = '  -2.28%  '
$ws.Range('D48').Value = '9.69'
$ws.Range('E48').Value = '  -0.36%  '
$ws.Range('E49').Value = '  -1.04%  '
$ws.Range('E50').Value = '  +26.15%  '
$ws.Range('D51').Value = '0.997'
$ws.Range('E51').Value = '  -0.18%  '
